$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.186238
$ws.Range("H2").Value = 0.558714
$ws.Range("I2").Value = 0.05023668284714279
$ws.Range("J2").Value = 0.05023668284714279
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.315861666666667
$ws.Range("N2").Value = 3.947585
$ws.Range("O2").Value = 0.2754050739440597
$ws.Range("P2").Value = 0.2754050739440597
$ws.Range("Q2").Value = 0.2450634450766667
$ws.Range("R2").Value = 2.20557100569
$ws.Range("S2").Value = 0.01383543735422164
$ws.Range("T2").Value = 0.01383543735422164
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.186238
$ws.Range("H3").Value = 0.558714
$ws.Range("I3").Value = 0.05023668284714279
$ws.Range("J3").Value = 0.05023668284714279
$ws.Range("O3").Value = 0.3040809095127364
$ws.Range("P3").Value = 0.3040809095127364
$ws.Range("Q3").Value = 0.2705800376153333
$ws.Range("R3").Value = 2.435220338538
$ws.Range("S3").Value = 0.01527601621106207
$ws.Range("T3").Value = 0.01527601621106207
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.186238
$ws.Range("H4").Value = 0.558714
$ws.Range("I4").Value = 0.05023668284714279
$ws.Range("J4").Value = 0.05023668284714279
$ws.Range("M4").Value = 2.009179666666667
$ws.Range("N4").Value = 6.027539
$ws.Range("O4").Value = 0.4205140165432039
$ws.Range("P4").Value = 0.4205140165432039
$ws.Range("Q4").Value = 0.3741856027606667
$ws.Range("R4").Value = 3.367670424846
$ws.Range("S4").Value = 0.02112522928185909
$ws.Range("T4").Value = 0.02112522928185909
# Row 5
$ws.Range("I5").Value = 0.659992587420158
$ws.Range("J5").Value = 0.6599925874201579
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.315861666666667
$ws.Range("N5").Value = 3.947585
$ws.Range("O5").Value = 0.2754050739440597
$ws.Range("P5").Value = 0.2754050739440597
$ws.Range("Q5").Value = 3.219560847406667
$ws.Range("R5").Value = 28.97604762666
$ws.Range("S5").Value = 0.1817653073409799
$ws.Range("T5").Value = 0.1817653073409799
# Row 6
$ws.Range("I6").Value = 0.659992587420158
$ws.Range("J6").Value = 0.6599925874201579
$ws.Range("O6").Value = 0.3040809095127364
$ws.Range("P6").Value = 0.3040809095127364
$ws.Range("S6").Value = 0.2006911462543858
$ws.Range("T6").Value = 0.2006911462543858
# Row 7
$ws.Range("I7").Value = 0.659992587420158
$ws.Range("J7").Value = 0.6599925874201579
$ws.Range("M7").Value = 2.009179666666667
$ws.Range("N7").Value = 6.027539
$ws.Range("O7").Value = 0.4205140165432039
$ws.Range("P7").Value = 0.4205140165432039
$ws.Range("Q7").Value = 4.915924184182668
$ws.Range("R7").Value = 44.24331765764401
$ws.Range("S7").Value = 0.2775361338247923
$ws.Range("T7").Value = 0.2775361338247922
# Row 8
$ws.Range("G8").Value = 1.074241333333333
$ws.Range("H8").Value = 3.222724
$ws.Range("I8").Value = 0.2897707297326994
$ws.Range("J8").Value = 0.2897707297326994
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.315861666666667
$ws.Range("N8").Value = 3.947585
$ws.Range("O8").Value = 0.2754050739440597
$ws.Range("P8").Value = 0.2754050739440597
$ws.Range("Q8").Value = 1.413552991282222
$ws.Range("R8").Value = 12.72197692154
$ws.Range("S8").Value = 0.07980432924885821
$ws.Range("T8").Value = 0.07980432924885821
# Row 9
$ws.Range("G9").Value = 1.074241333333333
$ws.Range("H9").Value = 3.222724
$ws.Range("I9").Value = 0.2897707297326994
$ws.Range("J9").Value = 0.2897707297326994
$ws.Range("O9").Value = 0.3040809095127364
$ws.Range("P9").Value = 0.3040809095127364
$ws.Range("Q9").Value = 1.560735512523111
$ws.Range("R9").Value = 14.046619612708
$ws.Range("S9").Value = 0.08811374704728855
$ws.Range("T9").Value = 0.08811374704728855
# Row 10
$ws.Range("G10").Value = 1.074241333333333
$ws.Range("H10").Value = 3.222724
$ws.Range("I10").Value = 0.2897707297326994
$ws.Range("J10").Value = 0.2897707297326994
$ws.Range("M10").Value = 2.009179666666667
$ws.Range("N10").Value = 6.027539
$ws.Range("O10").Value = 0.4205140165432039
$ws.Range("P10").Value = 0.4205140165432039
$ws.Range("Q10").Value = 2.158343844026223
$ws.Range("R10").Value = 19.425094596236
$ws.Range("S10").Value = 0.1218526534365526
$ws.Range("T10").Value = 0.1218526534365526
